$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D14: update raw value
$ws.Range("D14").Value = -45758211.060000002

# D18: replace hard-coded value with a SUM formula
$ws.Range("D18").Formula = "=SUM(D12:D17)"

# D21: replace hard-coded value with a SUM formula
$ws.Range("D21").Formula = "=SUM(D18:D20)"

# D22: update raw value
$ws.Range("D22").Value = -32201025

$excel.CalculateFull()
